# Increment the "想去人数" (F column) values by 1 for specific rows
# in both the "展览" and "全部类型" worksheets, matching the committed diff.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")
$rows = @(7, 12, 17, 26, 29)

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($r in $rows) {
        $cell = $ws.Cells.Item($r, 6)  # Column F = 6
        $cell.Value = $cell.Value2 + 1
    }
}
